$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Refresh rows 2-8 with the latest scrape (2025-09-11 06:26:24) ---
$ws.Range("A2").Value = '2025-09-11 06:26:24'
$ws.Range("B2").Value = '【日本人限定/継続案件】Node.jsエンジニア募集(スクレイピング機能開発)'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5391607'
$ws.Range("G2").Value = 155
$ws.Range("H2").Value = '◆開発,Node.js'

$ws.Range("A3").Value = '2025-09-11 06:26:24'
$ws.Range("B3").Value = '【急募】SharePoint+Power Platformでの不動産賃貸管理システム構築'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5391490'
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = '◇管理'

$ws.Range("A4").Value = '2025-09-11 06:26:24'
$ws.Range("B4").Value = '【急募】Salesforce・MA・CRMコンサルタント経験者を探しています!'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5371747'
$ws.Range("G4").Value = 48
$ws.Range("H4").Value = '◆コンサル'

$ws.Range("A5").Value = '2025-09-11 06:26:24'
$ws.Range("B5").Value = '初回 【フルリモート】フリーランスエンジニア募集'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5391489'
$ws.Range("G5").Value = 25
$ws.Range("H5").ClearContents()

$ws.Range("A6").Value = '2025-09-11 06:26:24'
$ws.Range("B6").Value = '要件定義や基本設計ができる方(1人月、約2年アサイン予定)'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5391221'
$ws.Range("G6").Value = 25
$ws.Range("H6").ClearContents()

$ws.Range("A7").Value = '2025-09-11 06:26:24'
$ws.Range("B7").Value = '【講師募集】Gensparkを使ったWEB構築チュートリアル募集'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5390165'
$ws.Range("G7").Value = 13
$ws.Range("H7").ClearContents()

$ws.Range("A8").Value = '2025-09-11 06:26:24'
$ws.Range("B8").Value = '【急募】Googleアナリティクス連携の専門家を探しています'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5391267'
$ws.Range("G8").Value = 10
$ws.Range("H8").ClearContents()

# --- Drop the stale listings that rolled off this scrape (old rows 9-22) ---
$ws.Range("A9:H22").EntireRow.Delete()

# --- Rebuild the hyperlinks on column F to match the new URLs ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5391607')
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5391490')
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5371747')
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5391489')
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5391221')
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5390165')
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5391267')

# --- Column width tweaks (title column narrower, skill column narrower) ---
$ws.Columns.Item(2).ColumnWidth = 45.16666666666667
$ws.Columns.Item(8).ColumnWidth = 12.16666666666667

